# GDP Growth Rates.xlsx - "files from canada june 22" update
#
# Semantic changes applied (per the supplied OOXML diff):
#   1. The "Data" sheet gains a new 2018 data column (G2:G4) used as the
#      base year for the November-GDP comparison table.
#   2. The workbook's active tab moves from "About" to "Data", and each
#      sheet's remembered selection is updated to reflect where the user
#      left off (About stays on B8; Data moves from C14 to C2;
#      GDPGR-alternate moves from E20 to B2; GDPGR-bau moves from C28 to B2).

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Data")
$wsAlt = $wb.Worksheets.Item("GDPGR-alternate")
$wsBau = $wb.Worksheets.Item("GDPGR-bau")

# --- New 2018 column of source data on the Data sheet ---
$wsData.Range("G2").Value = 2018
$wsData.Range("G3").Value = 1964629
$wsData.Range("G4").Value = 1922693

# --- Restore each sheet's remembered selection ---
$wsAbout.Range("B8").Select()
$wsAlt.Range("B2").Select()
$wsBau.Range("B2").Select()

# Activate "Data" last so it becomes the active/selected tab (matches
# workbookView activeTab + the tabSelected sheetView flag moving to Data).
$wsData.Activate()
$wsData.Range("C2").Select()
